# Applies the "475treatedDataset.xlsx" edit:
#  1. Rename the worksheet to reflect the filtered ("less2000") dataset.
#  2. Clear the cells that correspond to subjects/time-points whose tumour
#     volume exceeded 2000 (the "less2000" filter), in columns:
#       I/J  = CARBO18
#       M/N  = GEM25
#       Q/R  = GEM27
#       AC/AD = TRB40

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "475treated_less2000_GrowDataFil"

# Row 22: drop GEM25 (M22:N22) and GEM27 (Q22:R22)
$ws.Range("M22:N22").ClearContents()
$ws.Range("Q22:R22").ClearContents()

# Row 23: drop CARBO18 (I23:J23), GEM27 (Q23:R23), TRB40 (AC23:AD23)
$ws.Range("I23:J23").ClearContents()
$ws.Range("Q23:R23").ClearContents()
$ws.Range("AC23:AD23").ClearContents()

# Row 24: drop CARBO18 (I24:J24), TRB40 (AC24:AD24)
$ws.Range("I24:J24").ClearContents()
$ws.Range("AC24:AD24").ClearContents()

# Row 25: drop TRB40 (AC25:AD25)
$ws.Range("AC25:AD25").ClearContents()
